$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.664.80"
$ws.Range("E2").Value = "  -5.92%  "
$ws.Range("D3").Value = "1.805.58"
$ws.Range("E3").Value = "  -5.21%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "'274.88"
$ws.Range("E5").Value = "  -10.12%  "
$ws.Range("D6").Value = "'0.9996"
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("D7").Value = "'0.5045"
$ws.Range("E7").Value = "  -6.70%  "
$ws.Range("D8").Value = "'0.3506"
$ws.Range("E8").Value = "  -7.95%  "
$ws.Range("D9").Value = "'44.00"
$ws.Range("E9").Value = "  -3.94%  "
$ws.Range("D10").Value = "'0.06614"
$ws.Range("E10").Value = "  -9.45%  "
$ws.Range("D11").Value = "'19.87"
$ws.Range("E11").Value = "  -9.85%  "
$ws.Range("D12").Value = "'0.8325"
$ws.Range("E12").Value = "  -7.89%  "
$ws.Range("D13").Value = "'0.07793"
$ws.Range("E13").Value = "  -4.77%  "
$ws.Range("D14").Value = "1.799.32"
$ws.Range("E14").Value = "  +33.12%  "
$ws.Range("D15").Value = "'5.053"
$ws.Range("E15").Value = "  -5.49%  "
$ws.Range("D16").Value = "'87.31"
$ws.Range("E16").Value = "  -8.59%  "
$ws.Range("D17").Value = "'0.9995"
$ws.Range("E17").Value = "  +0.09%  "
$ws.Range("E18").Value = "  -6.47%  "
$ws.Range("D19").Value = "'0.9996"
$ws.Range("E19").Value = "  +0.09%  "
$ws.Range("D20").Value = "'0.000007960"
$ws.Range("E20").Value = "  -8.02%  "
$ws.Range("D21").Value = "25.723.85"
$ws.Range("E21").Value = "  -5.62%  "
$ws.Range("D22").Value = "'4.713"
$ws.Range("E22").Value = "  -6.66%  "
$ws.Range("D23").Value = "'9.965"
$ws.Range("E23").Value = "  -7.84%  "
$ws.Range("D24").Value = "'6.048"
$ws.Range("E24").Value = "  -7.11%  "
$ws.Range("E25").Value = "  -4.36%  "
$ws.Range("D26").Value = "'2.111"
$ws.Range("E26").Value = "  -8.40%  "
$ws.Range("D27").Value = "'1.651"
$ws.Range("E27").Value = "  -5.66%  "
$ws.Range("D28").Value = "'16.89"
$ws.Range("E28").Value = "  -7.90%  "
$ws.Range("D29").Value = "'108.30"
$ws.Range("E29").Value = "  -7.19%  "
$ws.Range("D30").Value = "'4.317"
$ws.Range("E30").Value = "  -10.90%  "
$ws.Range("D31").Value = "'4.186"
$ws.Range("E31").Value = "  -10.37%  "
$ws.Range("D32").Value = "'0.08773"
$ws.Range("E32").Value = "  -4.61%  "
$ws.Range("D33").Value = "'0.04786"
$ws.Range("E33").Value = "  -5.60%  "
$ws.Range("D34").Value = "'0.7235"
$ws.Range("E34").Value = "  -12.23%  "
$ws.Range("D35").Value = "'1.129"
$ws.Range("E35").Value = "  -7.81%  "
$ws.Range("D36").Value = "'2.871"
$ws.Range("E36").Value = "  -4.69%  "
$ws.Range("D37").Value = "'0.9988"
$ws.Range("E37").Value = "  +0.21%  "
$ws.Range("D38").Value = "'3.025"
$ws.Range("E38").Value = "  -8.92%  "
$ws.Range("D39").Value = "'0.01852"
$ws.Range("E39").Value = "  -7.47%  "
$ws.Range("D40").Value = "'0.5164"
$ws.Range("E40").Value = "  -14.46%  "
$ws.Range("D41").Value = "'2.273"
$ws.Range("E41").Value = "  -15.82%  "
$ws.Range("D42").Value = "'0.9436"
$ws.Range("E42").Value = "  -12.16%  "
$ws.Range("D43").Value = "'112.77"
$ws.Range("E43").Value = "  -2.81%  "
$ws.Range("E44").Value = "  -7.87%  "
$ws.Range("D45").Value = "'7.997"
$ws.Range("E45").Value = "  -13.84%  "
$ws.Range("D46").Value = "'0.9992"
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("D47").Value = "'0.4547"
$ws.Range("E47").Value = "  -12.14%  "
$ws.Range("D48").Value = "'0.1373"
$ws.Range("E48").Value = "  -10.43%  "
$ws.Range("D49").Value = "'9.255"
$ws.Range("E49").Value = "  -9.46%  "
$ws.Range("D50").Value = "'36.01"
$ws.Range("E50").Value = "  -5.54%  "
$ws.Range("D51").Value = "'1.487"
